# Weekly refresh of the "Mora" (blackberry) price sheet:
# a new, most-recent observation is inserted as row 12 (pushing the
# previously-existing rows 12-14 down to rows 13-15, unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12:14 down by one, creating a blank row 12 to fill in.
$ws.Rows("12:12").Insert()

# New weekly observation for row 12.
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44617
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100101
$ws.Range("H12").Value = "Berries"
$ws.Range("I12").Value = 100101008
$ws.Range("J12").Value = "Mora"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 90
$ws.Range("N12").Value = 6500
$ws.Range("O12").Value = 6500
$ws.Range("P12").Value = 6500
$ws.Range("Q12").Value = "`$/bandeja 2 kilos"
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 3250
$ws.Range("T12").Value = 2
